# Swap the last two comma-separated names in column G ("Recorded By")
# for every data row on the active sheet.
#
# Observed rule (derived from the target diff): wherever the "Recorded By"
# cell contains more than one comma-separated entry, the final two entries
# trade places while everything before them keeps its original order.
#   "System, dnasr281@gmail.com"                 -> "dnasr281@gmail.com, System"
#   "system, System, backup@backdoor.com"        -> "system, backup@backdoor.com, System"
#   "dnasr281@gmail.com, admin@admin.com"        -> "admin@admin.com, dnasr281@gmail.com"
# Cells with a single value are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $text = $cell.Value()

    if ($null -eq $text) { continue }
    if ($text -notlike "*,*") { continue }

    $parts = $text -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $last = $parts.Count - 1
    $tmp = $parts[$last]
    $parts[$last] = $parts[$last - 1]
    $parts[$last - 1] = $tmp

    $newText = [string]::Join(", ", $parts)
    $cell.Value = $newText
}
